$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New prediction columns N, P, R for rows 6-30 (plain-decimal literals; the
# interpreter's numeric literal grammar does not accept scientific notation)
$data = @{
    6  = @{ N = 0.086858; P = 0.089991; R = 0.090862 }
    7  = @{ N = 0.053406; P = 0.047836; R = 0.047637 }
    8  = @{ N = 0.132339; P = 0.126905; R = 0.129961 }
    9  = @{ N = 0.109937; P = 0.111815; R = 0.115898 }
    10 = @{ N = 0.066725; P = 0.075671; R = 0.072316 }
    11 = @{ N = 0.119032; P = 0.109262; R = 0.103513 }
    12 = @{ N = 0.066083; P = 0.066812; R = 0.066757 }
    13 = @{ N = 0.073587; P = 0.073042; R = 0.074964 }
    14 = @{ N = 0.007338; P = 0.012051; R = 0.015443 }
    15 = @{ N = 0.044155; P = 0.039912; R = 0.041798 }
    16 = @{ N = 0.047042; P = 0.044472; R = 0.046376 }
    17 = @{ N = 0.047726; P = 0.04861;  R = 0.044171 }
    18 = @{ N = 0.05891;  P = 0.05996;  R = 0.060598 }
    19 = @{ N = 0.040124; P = 0.03947;  R = 0.040069 }
    20 = @{ N = 0.050401; P = 0.045956; R = 0.047385 }
    21 = @{ N = 0.045982; P = 0.049167; R = 0.046159 }
    22 = @{ N = 0.048343; P = 0.05451;  R = 0.054074 }
    23 = @{ N = 0.093286; P = 0.092574; R = 0.0921 }
    24 = @{ N = 0.09827;  P = 0.09757;  R = 0.095497 }
    25 = @{ N = 0.040447; P = 0.041928; R = 0.041709 }
    26 = @{ N = 0.061883; P = 0.066809; R = 0.065483 }
    27 = @{ N = 0.01013;  P = 0.010944; R = 0.013062 }
    28 = @{ N = 0.04099;  P = 0.038518; R = 0.033201 }
    29 = @{ N = 0.031836; P = 0.034331; R = 0.030514 }
    30 = @{ N = 0.008547; P = -0.000795; R = 0.003149 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
}

# Error % formulas for rows 27-30, columns O, Q, S (mirrors existing C/E/G/I/K/M pattern)
foreach ($row in 27..30) {
    $ws.Range("O$row").Formula = "=ABS(N$row-`$A$row)/`$A$row"
    $ws.Range("Q$row").Formula = "=ABS(P$row-`$A$row)/`$A$row"
    $ws.Range("S$row").Formula = "=ABS(R$row-`$A$row)/`$A$row"
}

# Average row 31
$ws.Range("O31").Formula = "=AVERAGE(O27:O30)"
$ws.Range("Q31").Formula = "=AVERAGE(Q27:Q30)"
$ws.Range("S31").Formula = "=AVERAGE(S27:S30)"

# Match the percent number format already used by the other error columns
$ws.Range("O27:O31").NumberFormat = "0%"
$ws.Range("Q27:Q31").NumberFormat = "0%"
$ws.Range("S27:S31").NumberFormat = "0%"

$ws.Range("A1").Select()
